# Updates cryptos list values (prices & 1h volume changes) per the
# Jan 13 2024 GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.682.71"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -6.96%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.551.42"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -3.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.86"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -3.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.46"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -5.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -3.59%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -5.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.05"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.24%  "

# Row 11
$ws.Range("E11").Value = "  -5.00%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.73"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -4.34%  "

# Row 13
$ws.Range("E13").Value = "  +5.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.938.11"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -4.42%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.548.84"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -4.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.877"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -5.59%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.34"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -4.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.708.56"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -6.95%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -3.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.79"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.91"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -8.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -4.00%  "

# Row 25
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -4.75%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.22"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -5.88%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -4.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.17"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -3.41%  "

# Row 30
$ws.Range("E30").Value = "  -2.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.98"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -3.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.63"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.00%  "

# Row 33
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -7.28%  "

# Row 34
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.75"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -2.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.40"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -8.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0794"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -5.38%  "

# Row 37
$ws.Range("E37").Value = "  -6.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.42"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -4.74%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.119"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.92%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.94"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +6.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0311"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -4.82%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.88"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.15%  "

# Row 43
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.41"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.083.66"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "84.88"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -9.95%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.04"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.35%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.62"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +3.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.794.75"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -4.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.45"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -6.43%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.69"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.87%  "
